$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated voltage magnitude (vm_pu) results for case with 380 kV done
# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022192629141769
$ws.Range("D2").Value = 1.031106801433389
$ws.Range("E2").Value = 1.022973820019877
$ws.Range("F2").Value = 1.037112412676056
$ws.Range("I2").Value = 1.024935622675778
$ws.Range("J2").Value = 1.027379847277369
$ws.Range("K2").Value = 1.033915877480283
$ws.Range("L2").Value = 1.025806623784877
$ws.Range("M2").Value = 1.039904227029816
$ws.Range("N2").Value = 1.013190409921538
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023530341893124
$ws.Range("D3").Value = 1.032476280520569
$ws.Range("E3").Value = 1.024116507424092
$ws.Range("F3").Value = 1.038794557677048
$ws.Range("I3").Value = 1.025132670580908
$ws.Range("J3").Value = 1.02835348688739
$ws.Range("K3").Value = 1.035092346893526
$ws.Range("L3").Value = 1.026755174863784
$ws.Range("M3").Value = 1.041393803183955
$ws.Range("N3").Value = 1.013510772564229
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024389037690743
$ws.Range("D4").Value = 1.033352565983874
$ws.Range("E4").Value = 1.024850348983022
$ws.Range("F4").Value = 1.039864524779623
$ws.Range("I4").Value = 1.025249844857602
$ws.Range("J4").Value = 1.028976228645368
$ws.Range("K4").Value = 1.035843244259318
$ws.Range("L4").Value = 1.027362840909987
$ws.Range("M4").Value = 1.042338755614989
$ws.Range("N4").Value = 1.013715665731523
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024748402945398
$ws.Range("D5").Value = 1.033718620325351
$ws.Range("E5").Value = 1.025157542544762
$ws.Range("F5").Value = 1.040309951001723
$ws.Range("I5").Value = 1.025296640014843
$ws.Range("J5").Value = 1.029236306366369
$ws.Range("K5").Value = 1.03615646477105
$ws.Range("L5").Value = 1.027616856329176
$ws.Range("M5").Value = 1.042731524861771
$ws.Range("N5").Value = 1.013801233196304
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02480864690146
$ws.Range("D6").Value = 1.033779946148541
$ws.Range("E6").Value = 1.0252090451266
$ws.Range("F6").Value = 1.040384483809439
$ws.Range("I6").Value = 1.025304352792473
$ws.Range("J6").Value = 1.02927987400535
$ws.Range("K6").Value = 1.036208912443006
$ws.Range("L6").Value = 1.027659422189342
$ws.Range("M6").Value = 1.04279721034195
$ws.Range("N6").Value = 1.013815567104619
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024393845927135
$ws.Range("D7").Value = 1.033357466367222
$ws.Range("E7").Value = 1.024854458858047
$ws.Range("F7").Value = 1.039870493776747
$ws.Range("I7").Value = 1.025250479812058
$ws.Range("J7").Value = 1.028979710563991
$ws.Range("K7").Value = 1.035847439150586
$ws.Range("L7").Value = 1.027366240743488
$ws.Range("M7").Value = 1.042344021406739
$ws.Range("N7").Value = 1.013716811318993
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022646155185633
$ws.Range("D8").Value = 1.031571679743152
$ws.Range("E8").Value = 1.023361156488123
$ws.Range("F8").Value = 1.037684755112653
$ws.Range("I8").Value = 1.025004359667409
$ws.Range("J8").Value = 1.027710408898115
$ws.Range("K8").Value = 1.034315629477689
$ws.Range("L8").Value = 1.026128465696659
$ws.Range("M8").Value = 1.040411573120449
$ws.Range("N8").Value = 1.013299179010472
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019512751126477
$ws.Range("D9").Value = 1.028348222293584
$ws.Range("E9").Value = 1.020686460847832
$ws.Range("F9").Value = 1.033689666529277
$ws.Range("I9").Value = 1.024491147364305
$ws.Range("J9").Value = 1.025417248903161
$ws.Range("K9").Value = 1.031535964875113
$ws.Range("L9").Value = 1.023899816157457
$ws.Range("M9").Value = 1.036859799376705
$ws.Range("N9").Value = 1.012544583507632
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01738626277433
$ws.Range("D10").Value = 1.02614595743701
$ws.Range("E10").Value = 1.018873061936417
$ws.Range("F10").Value = 1.030926909139105
$ws.Range("I10").Value = 1.024094883258631
$ws.Range("J10").Value = 1.023849264308511
$ws.Range("K10").Value = 1.029627136787867
$ws.Range("L10").Value = 1.02238099023734
$ws.Range("M10").Value = 1.034390731731472
$ws.Range("N10").Value = 1.012028559341924
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016456224024171
$ws.Range("D11").Value = 1.025179291258633
$ws.Range("E11").Value = 1.018080383656476
$ws.Range("F11").Value = 1.029706343215077
$ws.Range("I11").Value = 1.023910287849645
$ws.Range("J11").Value = 1.023160711694576
$ws.Range("K11").Value = 1.028786978036289
$ws.Range("L11").Value = 1.021715214729545
$ws.Range("M11").Value = 1.033296936294683
$ws.Range("N11").Value = 1.011801943308702
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016109344898492
$ws.Range("D12").Value = 1.024818226169545
$ws.Range("E12").Value = 1.017784800890036
$ws.Range("F12").Value = 1.02924926023089
$ws.Range("I12").Value = 1.023839749727162
$ws.Range("J12").Value = 1.022903482341722
$ws.Range("K12").Value = 1.028472822259271
$ws.Range("L12").Value = 1.021466672755369
$ws.Range("M12").Value = 1.032886884954315
$ws.Range("N12").Value = 1.011717282197749
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016183816478498
$ws.Range("D13").Value = 1.024895767102429
$ws.Range("E13").Value = 1.017848256710779
$ws.Range("F13").Value = 1.029347475049018
$ws.Range("I13").Value = 1.02385496987307
$ws.Range("J13").Value = 1.022958725896506
$ws.Range("K13").Value = 1.028540304628843
$ws.Range("L13").Value = 1.021520042519129
$ws.Range("M13").Value = 1.032975013829086
$ws.Range("N13").Value = 1.011735464427006
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016427580065974
$ws.Range("D14").Value = 1.025149486580665
$ws.Range("E14").Value = 1.01805597423744
$ws.Range("F14").Value = 1.029668636666067
$ws.Range("I14").Value = 1.023904497457982
$ws.Range("J14").Value = 1.023139479191584
$ws.Range("K14").Value = 1.028761052541878
$ws.Range("L14").Value = 1.021694695680743
$ws.Range("M14").Value = 1.033263118529851
$ws.Range("N14").Value = 1.011794955159052
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016577581569838
$ws.Range("D15").Value = 1.025305545024378
$ws.Range("E15").Value = 1.018183803125525
$ws.Range("F15").Value = 1.029866021303698
$ws.Range("I15").Value = 1.023934751322284
$ws.Range("J15").Value = 1.023250651600134
$ws.Range("K15").Value = 1.028896785415034
$ws.Range("L15").Value = 1.021802139759071
$ws.Range("M15").Value = 1.033440128318758
$ws.Range("N15").Value = 1.011831544710145
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017447788789317
$ws.Range("D16").Value = 1.026209833261007
$ws.Range("E16").Value = 1.0189255099842
$ws.Range("F16").Value = 1.031007397311913
$ws.Range("I16").Value = 1.024106858755466
$ws.Range("J16").Value = 1.023894756622525
$ws.Range("K16").Value = 1.029682605244051
$ws.Range("L16").Value = 1.022425002571237
$ws.Range("M16").Value = 1.034462798444366
$ws.Range("N16").Value = 1.012043531473145
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017991147611956
$ws.Range("D17").Value = 1.026773543402717
$ws.Range("E17").Value = 1.019388747176962
$ws.Range("F17").Value = 1.031716809770873
$ws.Range("I17").Value = 1.024211322713639
$ws.Range("J17").Value = 1.024296196411985
$ws.Range("K17").Value = 1.030171856620046
$ws.Range("L17").Value = 1.022813518900808
$ws.Range("M17").Value = 1.035097645363481
$ws.Range("N17").Value = 1.012175649190447
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018307188787927
$ws.Range("D18").Value = 1.027101087543604
$ws.Range("E18").Value = 1.019658227149983
$ws.Range("F18").Value = 1.032128260482032
$ws.Range("I18").Value = 1.024271000645321
$ws.Range("J18").Value = 1.024529424630834
$ws.Range("K18").Value = 1.030455916774984
$ws.Range("L18").Value = 1.023039352743877
$ws.Range("M18").Value = 1.035465564514944
$ws.Range("N18").Value = 1.012252405589111
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018414800471183
$ws.Range("D19").Value = 1.027212559310127
$ws.Range("E19").Value = 1.019749991742069
$ws.Range("F19").Value = 1.032268159920545
$ws.Range("I19").Value = 1.024291137034495
$ws.Range("J19").Value = 1.02460879336313
$ws.Range("K19").Value = 1.030552552505006
$ws.Range("L19").Value = 1.023116224466216
$ws.Range("M19").Value = 1.035590614143298
$ws.Range("N19").Value = 1.012278525960825
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017932942713813
$ws.Range("D20").Value = 1.026713193078419
$ws.Range("E20").Value = 1.01933912070118
$ws.Range("F20").Value = 1.03164093876303
$ws.Range("I20").Value = 1.024200244555064
$ws.Range("J20").Value = 1.024253221572478
$ws.Range("K20").Value = 1.030119500531731
$ws.Range("L20").Value = 1.022771915746195
$ws.Range("M20").Value = 1.035029778515272
$ws.Range("N20").Value = 1.012161505885759
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016355837255693
$ws.Range("D21").Value = 1.025074828108427
$ws.Range("E21").Value = 1.017994838468746
$ws.Range("F21").Value = 1.029574165529401
$ws.Range("I21").Value = 1.023889967371667
$ws.Range("J21").Value = 1.023086292663301
$ws.Range("K21").Value = 1.028696105599694
$ws.Range("L21").Value = 1.021643299200488
$ws.Range("M21").Value = 1.033178383364721
$ws.Range("N21").Value = 1.01177745010657
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015356008737945
$ws.Range("D22").Value = 1.024033118204998
$ws.Range("E22").Value = 1.017142986606275
$ws.Range("F22").Value = 1.028253201560738
$ws.Range("I22").Value = 1.023683470078486
$ws.Range("J22").Value = 1.022344077169121
$ws.Range("K22").Value = 1.027789087159868
$ws.Range("L22").Value = 1.020926486833377
$ws.Range("M22").Value = 1.03199251065133
$ws.Range("N22").Value = 1.011533163380818
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015886828476575
$ws.Range("D23").Value = 1.02458646171701
$ws.Range("E23").Value = 1.017595208361077
$ws.Range("F23").Value = 1.028955529818225
$ws.Range("I23").Value = 1.023794025947107
$ws.Range("J23").Value = 1.022738356804889
$ws.Range("K23").Value = 1.028271071981189
$ws.Range("L23").Value = 1.021307174156177
$ws.Range("M23").Value = 1.03262325417477
$ws.Range("N23").Value = 1.011662934394434
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017959245753233
$ws.Range("D24").Value = 1.026740466679049
$ws.Range("E24").Value = 1.01936154698837
$ws.Range("F24").Value = 1.031675228825562
$ws.Range("I24").Value = 1.024205254172568
$ws.Range("J24").Value = 1.024272642909068
$ws.Range("K24").Value = 1.030143162046429
$ws.Range("L24").Value = 1.022790716833681
$ws.Range("M24").Value = 1.035060451962476
$ws.Range("N24").Value = 1.012167897580347
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020329321364168
$ws.Range("D25").Value = 1.029190822129955
$ws.Range("E25").Value = 1.021383180087203
$ws.Range("F25").Value = 1.034739785244915
$ws.Range("I25").Value = 1.024633311541331
$ws.Range("J25").Value = 1.026016901901067
$ws.Range("K25").Value = 1.032264268079562
$ws.Range("L25").Value = 1.024481717471827
$ws.Range("M25").Value = 1.037795650214143
$ws.Range("N25").Value = 1.012741917646915
